# "added resolve count to affect idle time"
# - IDLE TIME sheet: add two new employees (MAKEDA.OLLIVIERRE, RARG046N.YEBOAH)
#   in alphabetical position, add a trailing "~" bucket row, and refresh the
#   idle-time totals for PATR5027.AMEH and ZAHIDGUL.MINHAS.
# - Total Units picked by hour sheet: the hourly pick counts were being stored
#   as negative numbers; flip them back to positive now that resolve counts
#   are factored into the idle-time calculation.

$wb = $excel.ActiveWorkbook

$idle = $wb.Worksheets.Item("IDLE TIME")

# Insert MAKEDA.OLLIVIERRE right before MARI882N.ABDELKADER (currently row 15)
$idle.Rows.Item(15).Insert()
$idle.Cells.Item(15, 1).Value = "MAKEDA.OLLIVIERRE"
$idle.Cells.Item(15, 2).Value = 149

# Update PATR5027.AMEH's total idle time (now shifted down to row 20)
$idle.Cells.Item(20, 2).Value = 138

# Insert RARG046N.YEBOAH right after PRINCE.FORSON (currently row 21)
$idle.Rows.Item(22).Insert()
$idle.Cells.Item(22, 1).Value = "RARG046N.YEBOAH"
$idle.Cells.Item(22, 2).Value = 159

# Update ZAHIDGUL.MINHAS's total idle time (now shifted down to row 31)
$idle.Cells.Item(31, 2).Value = 120

# Append the "~" catch-all bucket as the new last row
$idle.Cells.Item(33, 1).Value = "~"
$idle.Cells.Item(33, 2).Value = 34

# Total Units picked by hour: un-negate the hourly pick counts
$hours = $wb.Worksheets.Item("Total Units picked by hour")
$rng = $hours.Range("B2:E6")
for ($r = 1; $r -le $rng.Rows.Count; $r++) {
    for ($c = 1; $c -le $rng.Columns.Count; $c++) {
        $cell = $rng.Cells.Item($r, $c)
        $cell.Value = -1 * $cell.Value()
    }
}
